# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
#
# 1) IsabellaJourney (sheet1): soften/diversify a handful of PPV-edge lines.
# 2) Split the old "dickpic" sheet content: the dick-pic-reaction lines now
#    get their own fresh "dickpic" tab (untouched), while a *new* second
#    control tab ("cumcontrol2") is grown out of a copy of it with new
#    delay/sync/edge style copy. The original "cumcontrol" tab becomes
#    "cumcontrol1" and gets its own line-level rewrite.
# Final tab order: ... done1, done2, cumcontrol1, cumcontrol2, dickpic, boosters

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) IsabellaJourney text tweaks
# ---------------------------------------------------------------------------
$journey = $wb.Worksheets.Item("IsabellaJourney")

$journey.Range("B4").Value  = "I want us to finish together"
$journey.Range("B5").Value  = "wait"
$journey.Range("B8").Value  = "god"
$journey.Range("B10").Value = "one second"
$journey.Range("B11").Value = "I'm done holding back"
$journey.Range("B20").Value = "did you see it?"
$journey.Range("B22").Value = "one second"

# ---------------------------------------------------------------------------
# 2) Rename "cumcontrol" -> "cumcontrol1" and rewrite its copy lines
# ---------------------------------------------------------------------------
$cumcontrol1 = $wb.Worksheets.Item("cumcontrol")
$cumcontrol1.Name = "cumcontrol1"

$cumcontrol1.Range("B2").Value = "just hold on a little more, I want the last thing you see to be this"

$cumcontrol1.Range("B3").Value = "wait for me... I have one more thing and I want you to see it before we finish"
$cumcontrol1.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol1.Range("B4").Value = "stay with me love, I'm almost there too... watch this"
$cumcontrol1.Range("C4").Value = "SYNC variant. Send PPV."

$cumcontrol1.Range("B5").Value = "I want us to finish together... open this and let go with me"
$cumcontrol1.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol1.Range("B6").Value = "please don't finish yet... I'm not ready for this to be over"

$cumcontrol1.Range("B7").Value = "not yet love... I want this to last a little longer with you"
$cumcontrol1.Range("C7").Value = "CONTROL."

# ---------------------------------------------------------------------------
# 3) Duplicate "dickpic" right after "cumcontrol1" -> becomes "cumcontrol2",
#    rewritten with the delay/sync/edge control copy. The original "dickpic"
#    sheet is left completely untouched.
# ---------------------------------------------------------------------------
$dickpic = $wb.Worksheets.Item("dickpic")
$dickpic.Copy($null, $cumcontrol1)

$cumcontrol2 = $wb.Worksheets.Item("dickpic (2)")
$cumcontrol2.Name = "cumcontrol2"

$cumcontrol2.Range("A2").Value = "delay2"
$cumcontrol2.Range("B2").Value = "just a little longer for me love? the next one is special"
$cumcontrol2.Range("C2").Value = "DELAY variant."

$cumcontrol2.Range("A3").Value = "delay1"
$cumcontrol2.Range("B3").Value = "please wait... what I'm about to send, I want you to really take it in"
$cumcontrol2.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol2.Range("A4").Value = "sync2"
$cumcontrol2.Range("B4").Value = "I need you to see this before we both let go"
$cumcontrol2.Range("C4").Value = "SYNC variant."

$cumcontrol2.Range("A5").Value = "sync1"
$cumcontrol2.Range("B5").Value = "okay... together, right now... open this"
$cumcontrol2.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol2.Range("A6").Value = "edge2"
$cumcontrol2.Range("B6").Value = "don't rush... this is too good to end yet"
$cumcontrol2.Range("C6").Value = "EDGE variant."

$cumcontrol2.Range("A7").Value = "edge1"
$cumcontrol2.Range("B7").Value = "slow down love... I want to feel every second of this with you"
$cumcontrol2.Range("C7").Value = "CONTROL."

# ---------------------------------------------------------------------------
# 4) Make sure "dickpic" (the untouched original) sits right after
#    "cumcontrol2" so the final tab order is:
#    ... done1, done2, cumcontrol1, cumcontrol2, dickpic, boosters
# ---------------------------------------------------------------------------
$dickpic.Move($null, $cumcontrol2)
